$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,4).Value2 = '61.742.91'
$ws.Cells.Item(2,5).Value2 = '  -8.04%  '
$ws.Cells.Item(3,4).Value2 = '2.905.45'
$ws.Cells.Item(3,5).Value2 = '  -10.16%  '
$c = $ws.Cells.Item(4,4)
$c.NumberFormat = "@"
$c.Value = '1.01'
$c.ClearFormats()
$ws.Cells.Item(4,5).Value2 = '  +0.65%  '
$c = $ws.Cells.Item(5,4)
$c.NumberFormat = "@"
$c.Value = '521.68'
$c.ClearFormats()
$ws.Cells.Item(5,5).Value2 = '  -11.93%  '
$c = $ws.Cells.Item(6,4)
$c.NumberFormat = "@"
$c.Value = '124.41'
$c.ClearFormats()
$ws.Cells.Item(6,5).Value2 = '  -18.55%  '
$c = $ws.Cells.Item(7,4)
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.ClearFormats()
$ws.Cells.Item(7,5).Value2 = '  +0.25%  '
$ws.Cells.Item(8,4).Value2 = '2.898.82'
$ws.Cells.Item(8,5).Value2 = '  -10.19%  '
$c = $ws.Cells.Item(9,4)
$c.NumberFormat = "@"
$c.Value = '0.437'
$c.ClearFormats()
$ws.Cells.Item(9,5).Value2 = '  -19.67%  '
$c = $ws.Cells.Item(10,4)
$c.NumberFormat = "@"
$c.Value = '0.139'
$c.ClearFormats()
$ws.Cells.Item(10,5).Value2 = '  -19.14%  '
$c = $ws.Cells.Item(11,4)
$c.NumberFormat = "@"
$c.Value = '5.68'
$c.ClearFormats()
$ws.Cells.Item(11,5).Value2 = '  -12.79%  '
$c = $ws.Cells.Item(12,4)
$c.NumberFormat = "@"
$c.Value = '0.417'
$c.ClearFormats()
$ws.Cells.Item(12,5).Value2 = '  -15.60%  '
$c = $ws.Cells.Item(13,4)
$c.NumberFormat = "@"
$c.Value = '0.0000197'
$c.ClearFormats()
$ws.Cells.Item(13,5).Value2 = '  -19.41%  '
$c = $ws.Cells.Item(14,4)
$c.NumberFormat = "@"
$c.Value = '30.65'
$c.ClearFormats()
$ws.Cells.Item(14,5).Value2 = '  -21.39%  '
$ws.Cells.Item(15,4).Value2 = '3.428.07'
$ws.Cells.Item(15,5).Value2 = '  -9.07%  '
$ws.Cells.Item(16,4).Value2 = '62.100.09'
$ws.Cells.Item(16,5).Value2 = '  -7.62%  '
$c = $ws.Cells.Item(17,4)
$c.NumberFormat = "@"
$c.Value = '0.109'
$c.ClearFormats()
$ws.Cells.Item(17,5).Value2 = '  -4.59%  '
$ws.Cells.Item(18,4).Value2 = '2.950.57'
$ws.Cells.Item(18,5).Value2 = '  -8.96%  '
$c = $ws.Cells.Item(19,4)
$c.NumberFormat = "@"
$c.Value = '460.86'
$c.ClearFormats()
$ws.Cells.Item(19,5).Value2 = '  -13.41%  '
$c = $ws.Cells.Item(20,4)
$c.NumberFormat = "@"
$c.Value = '5.93'
$c.ClearFormats()
$ws.Cells.Item(20,5).Value2 = '  -16.81%  '
$c = $ws.Cells.Item(21,4)
$c.NumberFormat = "@"
$c.Value = '12.25'
$c.ClearFormats()
$ws.Cells.Item(21,5).Value2 = '  -17.86%  '
$c = $ws.Cells.Item(22,4)
$c.NumberFormat = "@"
$c.Value = '0.612'
$c.ClearFormats()
$ws.Cells.Item(22,5).Value2 = '  -19.23%  '
$c = $ws.Cells.Item(23,4)
$c.NumberFormat = "@"
$c.Value = '6.26'
$c.ClearFormats()
$ws.Cells.Item(23,5).Value2 = '  -21.29%  '
$c = $ws.Cells.Item(24,4)
$c.NumberFormat = "@"
$c.Value = '72.85'
$c.ClearFormats()
$ws.Cells.Item(24,5).Value2 = '  -15.03%  '
$ws.Cells.Item(25,2).Value2 = 'Dai'
$ws.Cells.Item(25,3).Value2 = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$c = $ws.Cells.Item(25,4)
$c.NumberFormat = "@"
$c.Value = '0.998'
$c.ClearFormats()
$ws.Cells.Item(25,5).Value2 = '  -0.32%  '
$ws.Cells.Item(26,2).Value2 = 'InternetComputer(DFINITY)'
$ws.Cells.Item(26,3).Value2 = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$c = $ws.Cells.Item(26,4)
$c.NumberFormat = "@"
$c.Value = '11.41'
$c.ClearFormats()
$ws.Cells.Item(26,5).Value2 = '  -17.51%  '
$c = $ws.Cells.Item(27,4)
$c.NumberFormat = "@"
$c.Value = '2.58'
$c.ClearFormats()
$ws.Cells.Item(27,5).Value2 = '  -19.16%  '
$c = $ws.Cells.Item(28,4)
$c.NumberFormat = "@"
$c.Value = '6.69'
$c.ClearFormats()
$ws.Cells.Item(28,5).Value2 = '  -17.22%  '
$c = $ws.Cells.Item(29,4)
$c.NumberFormat = "@"
$c.Value = '1.78'
$c.ClearFormats()
$ws.Cells.Item(29,5).Value2 = '  -18.53%  '
$c = $ws.Cells.Item(30,4)
$c.NumberFormat = "@"
$c.Value = '23.82'
$c.ClearFormats()
$ws.Cells.Item(30,5).Value2 = '  -18.21%  '
$c = $ws.Cells.Item(31,4)
$c.NumberFormat = "@"
$c.Value = '1.06'
$c.ClearFormats()
$ws.Cells.Item(31,5).Value2 = '  -8.40%  '
$c = $ws.Cells.Item(32,4)
$c.NumberFormat = "@"
$c.Value = '1.01'
$c.ClearFormats()
$ws.Cells.Item(32,5).Value2 = '  +0.29%  '
$c = $ws.Cells.Item(33,4)
$c.NumberFormat = "@"
$c.Value = '2.20'
$c.ClearFormats()
$ws.Cells.Item(33,5).Value2 = '  -17.46%  '
$c = $ws.Cells.Item(34,4)
$c.NumberFormat = "@"
$c.Value = '50.84'
$c.ClearFormats()
$ws.Cells.Item(34,5).Value2 = '  -5.21%  '
$c = $ws.Cells.Item(35,4)
$c.NumberFormat = "@"
$c.Value = '461.46'
$c.ClearFormats()
$ws.Cells.Item(35,5).Value2 = '  -15.09%  '
$c = $ws.Cells.Item(36,4)
$c.NumberFormat = "@"
$c.Value = '5.25'
$c.ClearFormats()
$ws.Cells.Item(36,5).Value2 = '  -18.47%  '
$c = $ws.Cells.Item(37,4)
$c.NumberFormat = "@"
$c.Value = '4.52'
$c.ClearFormats()
$ws.Cells.Item(37,5).Value2 = '  -21.38%  '
$c = $ws.Cells.Item(38,4)
$c.NumberFormat = "@"
$c.Value = '0.0373'
$c.ClearFormats()
$ws.Cells.Item(38,5).Value2 = '  -12.49%  '
$c = $ws.Cells.Item(39,4)
$c.NumberFormat = "@"
$c.Value = '0.0731'
$c.ClearFormats()
$ws.Cells.Item(39,5).Value2 = '  -15.02%  '
$c = $ws.Cells.Item(40,4)
$c.NumberFormat = "@"
$c.Value = '0.108'
$c.ClearFormats()
$ws.Cells.Item(40,5).Value2 = '  -12.74%  '
$c = $ws.Cells.Item(41,4)
$c.NumberFormat = "@"
$c.Value = '7.49'
$c.ClearFormats()
$ws.Cells.Item(41,5).Value2 = '  -19.35%  '
$ws.Cells.Item(42,4).Value2 = '2.591.01'
$ws.Cells.Item(42,5).Value2 = '  -11.77%  '
$ws.Cells.Item(43,5).Value2 = '  -0.20%  '
$c = $ws.Cells.Item(44,4)
$c.NumberFormat = "@"
$c.Value = '2.14'
$c.ClearFormats()
$ws.Cells.Item(44,5).Value2 = '  -19.87%  '
$ws.Cells.Item(45,2).Value2 = 'Monero'
$ws.Cells.Item(45,3).Value2 = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$c = $ws.Cells.Item(45,4)
$c.NumberFormat = "@"
$c.Value = '111.37'
$c.ClearFormats()
$ws.Cells.Item(45,5).Value2 = '  -5.83%  '
$ws.Cells.Item(46,2).Value2 = 'TheGraph'
$ws.Cells.Item(46,3).Value2 = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$c = $ws.Cells.Item(46,4)
$c.NumberFormat = "@"
$c.Value = '0.215'
$c.ClearFormats()
$ws.Cells.Item(46,5).Value2 = '  -18.63%  '
$ws.Cells.Item(47,2).Value2 = 'Stellar'
$ws.Cells.Item(47,3).Value2 = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$c = $ws.Cells.Item(47,4)
$c.NumberFormat = "@"
$c.Value = '0.0969'
$c.ClearFormats()
$ws.Cells.Item(47,5).Value2 = '  -15.62%  '
$ws.Cells.Item(48,2).Value2 = 'Fetch.AI'
$ws.Cells.Item(48,3).Value2 = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$c = $ws.Cells.Item(48,4)
$c.NumberFormat = "@"
$c.Value = '1.74'
$c.ClearFormats()
$ws.Cells.Item(48,5).Value2 = '  -19.21%  '
$ws.Cells.Item(49,2).Value2 = 'BitgetToken'
$ws.Cells.Item(49,3).Value2 = 'https://coinranking.com/coin/q7gMmMdLb+bitgettoken-bgb'
$c = $ws.Cells.Item(49,4)
$c.NumberFormat = "@"
$c.Value = '1.18'
$c.ClearFormats()
$ws.Cells.Item(49,5).Value2 = '  -5.53%  '
$ws.Cells.Item(50,2).Value2 = 'InjectiveProtocol'
$ws.Cells.Item(50,3).Value2 = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$c = $ws.Cells.Item(50,4)
$c.NumberFormat = "@"
$c.Value = '20.94'
$c.ClearFormats()
$ws.Cells.Item(50,5).Value2 = '  -21.14%  '
$ws.Cells.Item(51,2).Value2 = 'PEPE'
$ws.Cells.Item(51,3).Value2 = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Cells.Item(51,4).Value2 = '0.0₃0446'
$ws.Cells.Item(51,5).Value2 = '  -23.87%  '
